$d = $word.ActiveDocument

# Locate the first occurrence of the "2634:WRKDPSTWK" run (its own paragraph,
# which currently lives inside a paragraph whose paragraph-mark formatting is
# Arial/sz20). We replace that whole paragraph (the Find range plus the
# trailing paragraph mark) with three paragraphs:
#   1) "2634:WRKDPSTWK" promoted into its own paragraph, paragraph mark now
#      matching the run's own (Menlo/sz22) formatting.
#   2) a brand new paragraph "2577:AGE2" (Arial/sz22, with the tab-stop /
#      autoSpace paragraph formatting used elsewhere in this document).
#   3) the original (now emptied) paragraph, keeping its original Arial/sz20
#      paragraph-mark formatting but with no run left inside it.

$findRange = $d.Content
$found = $findRange.Find.Execute("2634:WRKDPSTWK", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target text '2634:WRKDPSTWK'"
}

# Extend the range to also cover the trailing paragraph mark so the whole
# paragraph (mark included) gets replaced by the XML below.
$target = $d.Range($findRange.Start, $findRange.End + 1)

$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$xml = '<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" ' + $w + '>' `
  + '<w:p>' `
    + '<w:pPr>' `
      + '<w:rPr>' `
        + '<w:rFonts w:ascii="Menlo" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Menlo" w:cs="Menlo"/>' `
        + '<w:color w:val="5982DB" w:themeColor="accent6"/>' `
        + '<w:sz w:val="22"/>' `
        + '<w:szCs w:val="22"/>' `
      + '</w:rPr>' `
    + '</w:pPr>' `
    + '<w:r>' `
      + '<w:rPr>' `
        + '<w:rFonts w:ascii="Menlo" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Menlo" w:cs="Menlo"/>' `
        + '<w:color w:val="5982DB" w:themeColor="accent6"/>' `
        + '<w:sz w:val="22"/>' `
        + '<w:szCs w:val="22"/>' `
      + '</w:rPr>' `
      + '<w:t>2634:WRKDPSTWK</w:t>' `
    + '</w:r>' `
  + '</w:p>' `
  + '<w:p>' `
    + '<w:pPr>' `
      + '<w:tabs>' `
        + '<w:tab w:val="left" w:pos="560"/>' `
        + '<w:tab w:val="left" w:pos="1120"/>' `
        + '<w:tab w:val="left" w:pos="1680"/>' `
        + '<w:tab w:val="left" w:pos="2240"/>' `
        + '<w:tab w:val="left" w:pos="2800"/>' `
        + '<w:tab w:val="left" w:pos="3360"/>' `
        + '<w:tab w:val="left" w:pos="3920"/>' `
        + '<w:tab w:val="left" w:pos="4480"/>' `
        + '<w:tab w:val="left" w:pos="5040"/>' `
        + '<w:tab w:val="left" w:pos="5600"/>' `
        + '<w:tab w:val="left" w:pos="6160"/>' `
        + '<w:tab w:val="left" w:pos="6720"/>' `
      + '</w:tabs>' `
      + '<w:autoSpaceDE w:val="0"/>' `
      + '<w:autoSpaceDN w:val="0"/>' `
      + '<w:adjustRightInd w:val="0"/>' `
      + '<w:rPr>' `
        + '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>' `
        + '<w:color w:val="5982DB" w:themeColor="accent6"/>' `
        + '<w:sz w:val="22"/>' `
        + '<w:szCs w:val="22"/>' `
      + '</w:rPr>' `
    + '</w:pPr>' `
    + '<w:r>' `
      + '<w:rPr>' `
        + '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>' `
        + '<w:color w:val="5982DB" w:themeColor="accent6"/>' `
        + '<w:sz w:val="22"/>' `
        + '<w:szCs w:val="22"/>' `
      + '</w:rPr>' `
      + '<w:t>2577:AGE2</w:t>' `
    + '</w:r>' `
  + '</w:p>' `
  + '<w:p>' `
    + '<w:pPr>' `
      + '<w:rPr>' `
        + '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>' `
        + '<w:color w:val="5982DB" w:themeColor="accent6"/>' `
        + '<w:sz w:val="20"/>' `
        + '<w:szCs w:val="20"/>' `
      + '</w:rPr>' `
    + '</w:pPr>' `
  + '</w:p>' `
  + '</pkg:xmlData>'

$target.InsertXML($xml) | Out-Null
